$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated odds values per the 2024-10-12 FlashScore refresh.
# Cell -> new value, grouped by row for readability.

# Row 11
$ws.Range("J11").Value = 4.3
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 2.57
$ws.Range("W11").Value = 9.25
$ws.Range("X11").Value = 20
$ws.Range("AA11").Value = 40
$ws.Range("AE11").Value = 17.5
$ws.Range("AF11").Value = 100
$ws.Range("AG11").Value = 5.8
$ws.Range("AK11").Value = 18
$ws.Range("AL11").Value = 37
$ws.Range("AP11").Value = 30
$ws.Range("AR11").Value = 175
$ws.Range("AS11").Value = 450
$ws.Range("AU11").Value = 7.6
$ws.Range("AX11").Value = 10
$ws.Range("AZ11").Value = 40

# Row 12
$ws.Range("H12").Value = 3.7
$ws.Range("J12").Value = 2.2
$ws.Range("K12").Value = 2.22
$ws.Range("L12").Value = 4.65
$ws.Range("N12").Value = 11.5
$ws.Range("U12").Value = 1.65
$ws.Range("V12").Value = 2
$ws.Range("X12").Value = 8.75
$ws.Range("Z12").Value = 13.5
$ws.Range("AC12").Value = 12
$ws.Range("AD12").Value = 7.3
$ws.Range("AE12").Value = 14
$ws.Range("AG12").Value = 14
$ws.Range("AI12").Value = 14.5
$ws.Range("AN12").Value = 3.6
$ws.Range("AO12").Value = 8
$ws.Range("AP12").Value = 16
$ws.Range("AQ12").Value = 26
$ws.Range("AT12").Value = 2.92
$ws.Range("AU12").Value = 7
$ws.Range("AV12").Value = 60
$ws.Range("AW12").Value = 6.3
$ws.Range("AX12").Value = 25

# Row 15
$ws.Range("G15").Value = 1.62
$ws.Range("H15").Value = 3.6
$ws.Range("I15").Value = 5
$ws.Range("K15").Value = 2.2
$ws.Range("M15").Value = 1.07
$ws.Range("N15").Value = 9
$ws.Range("Q15").Value = 2.05
$ws.Range("R15").Value = 1.75
$ws.Range("X15").Value = 7
$ws.Range("AD15").Value = 7.5
$ws.Range("AN15").Value = 3.5
$ws.Range("AO15").Value = 8.5
$ws.Range("AW15").Value = 7
$ws.Range("BA15").Value = 151

# Row 35
$ws.Range("I35").Value = 4.25
$ws.Range("J35").Value = 2.3
$ws.Range("K35").Value = 2.15
$ws.Range("L35").Value = 4.55
$ws.Range("Q35").Value = 1.95
$ws.Range("R35").Value = 1.75
$ws.Range("W35").Value = 6.6
$ws.Range("X35").Value = 8.25
$ws.Range("AA35").Value = 14.5
$ws.Range("AD35").Value = 6.7
$ws.Range("AE35").Value = 16
$ws.Range("AF35").Value = 80
$ws.Range("AG35").Value = 11.25
$ws.Range("AH35").Value = 23
$ws.Range("AK35").Value = 45
$ws.Range("AL35").Value = 50
$ws.Range("AM35").Value = 700
$ws.Range("AO35").Value = 8.5
$ws.Range("AP35").Value = 17
$ws.Range("AQ35").Value = 29
$ws.Range("AR35").Value = 55
$ws.Range("AS35").Value = 200
$ws.Range("AT35").Value = 2.62
$ws.Range("AU35").Value = 7.2
$ws.Range("AV35").Value = 65
$ws.Range("AX35").Value = 24
$ws.Range("AY35").Value = 29
$ws.Range("AZ35").Value = 150
$ws.Range("BA35").Value = 175

# Row 36
$ws.Range("G36").Value = 2.9
$ws.Range("H36").Value = 3.15
$ws.Range("I36").Value = 2.35
$ws.Range("J36").Value = 3.45
$ws.Range("L36").Value = 2.95
$ws.Range("W36").Value = 9.5
$ws.Range("X36").Value = 15.5
$ws.Range("AA36").Value = 24
$ws.Range("AB36").Value = 30
$ws.Range("AG36").Value = 8.5
$ws.Range("AK36").Value = 19
$ws.Range("AL36").Value = 26
$ws.Range("AO36").Value = 16
$ws.Range("AP36").Value = 22
$ws.Range("AQ36").Value = 75
$ws.Range("AX36").Value = 12.5
$ws.Range("AY36").Value = 19.5
$ws.Range("AZ36").Value = 50
$ws.Range("BA36").Value = 80

# Row 37
$ws.Range("K37").Value = 2.02
$ws.Range("AA37").Value = 19.5
$ws.Range("AH37").Value = 14.5
$ws.Range("AT37").Value = 2.4
$ws.Range("AU37").Value = 7.5
$ws.Range("BB37").Value = 400
